# Update portfolio example: change purchase dates and avg purchase prices
# for the 5 holdings (rows 2-6). Dependent formulas (PURCHASE COST,
# CURRENT VALUE, UNREALIZED P&L) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "DATE PURCHASED" value (1/2/2020) for every holding row
$newDate = 43832

$ws.Range("B2").Value2 = $newDate
$ws.Range("B3").Value2 = $newDate
$ws.Range("B4").Value2 = $newDate
$ws.Range("B5").Value2 = $newDate
$ws.Range("B6").Value2 = $newDate

# New "AVG PURCHASE PRICE" values
$ws.Range("D2").Value2 = 29.37
$ws.Range("D3").Value2 = 74.06
$ws.Range("D4").Value2 = 158.78
$ws.Range("D5").Value2 = 93.75
$ws.Range("D6").Value2 = 145.87

# D3 and D5 swap from the "Yahoo Sans Finance" font to plain Arial,
# while keeping the same fill/alignment; F5 keeps the original font.
$ws.Range("D3").Font.Name = "Arial"
$ws.Range("D5").Font.Name = "Arial"
